$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "37.273.88"
Set-TextValue $ws.Cells.Item(2, 5) "  +3.50%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "2.036.86"
Set-TextValue $ws.Cells.Item(3, 5) "  +1.03%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  -0.24%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "253.29"
Set-TextValue $ws.Cells.Item(5, 5) "  +4.70%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "0.644"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.01%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "63.01"
Set-TextValue $ws.Cells.Item(7, 5) "  +15.38%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 5) "  +0.07%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "59.30"
Set-TextValue $ws.Cells.Item(9, 5) "  +1.76%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.375"
Set-TextValue $ws.Cells.Item(10, 5) "  +4.39%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.0750"
Set-TextValue $ws.Cells.Item(11, 5) "  +2.45%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.104"
Set-TextValue $ws.Cells.Item(12, 5) "  -0.83%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.904"
Set-TextValue $ws.Cells.Item(13, 5) "  +1.76%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "15.08"
Set-TextValue $ws.Cells.Item(14, 5) "  +7.13%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "2.322.30"
Set-TextValue $ws.Cells.Item(15, 5) "  +0.41%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "20.53"
Set-TextValue $ws.Cells.Item(16, 5) "  +19.64%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "5.50"
Set-TextValue $ws.Cells.Item(17, 5) "  +4.61%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "2.029.13"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.70%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "37.143.66"
Set-TextValue $ws.Cells.Item(19, 5) "  +3.53%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "72.63"
Set-TextValue $ws.Cells.Item(20, 5) "  +2.41%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "0.0₃0870"
Set-TextValue $ws.Cells.Item(21, 5) "  +2.87%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "5.33"
Set-TextValue $ws.Cells.Item(22, 5) "  +4.16%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "235.83"
Set-TextValue $ws.Cells.Item(23, 5) "  +0.07%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "2.80"
Set-TextValue $ws.Cells.Item(24, 5) "  +24.41%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 5) "  -0.25%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "2.34"
Set-TextValue $ws.Cells.Item(26, 5) "  +0.01%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "9.59"
Set-TextValue $ws.Cells.Item(27, 5) "  +4.77%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "164.93"
Set-TextValue $ws.Cells.Item(28, 5) "  +1.05%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "19.78"
Set-TextValue $ws.Cells.Item(29, 5) "  +0.32%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "0.121"
Set-TextValue $ws.Cells.Item(30, 5) "  +1.28%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "5.18"
Set-TextValue $ws.Cells.Item(31, 5) "  +6.49%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "1.20"
Set-TextValue $ws.Cells.Item(32, 5) "  +4.83%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "0.110"
Set-TextValue $ws.Cells.Item(33, 5) "  +23.59%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "4.72"
Set-TextValue $ws.Cells.Item(34, 5) "  +9.98%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "0.0613"
Set-TextValue $ws.Cells.Item(35, 5) "  +3.45%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "2.49"
Set-TextValue $ws.Cells.Item(36, 5) "  +14.86%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 5) "  -0.12%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 2) "THORChain"
Set-TextValue $ws.Cells.Item(38, 3) "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Cells.Item(38, 4) "6.01"
Set-TextValue $ws.Cells.Item(38, 5) "  +22.13%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 2) "WEMIXToken"
Set-TextValue $ws.Cells.Item(39, 3) "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Cells.Item(39, 4) "1.81"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.79%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.105"
Set-TextValue $ws.Cells.Item(40, 5) "  +17.79%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "1.24"
Set-TextValue $ws.Cells.Item(41, 5) "  +4.22%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "2.79"
Set-TextValue $ws.Cells.Item(42, 5) "  +25.00%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 5) "  +1.49%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "1.14"
Set-TextValue $ws.Cells.Item(44, 5) "  +4.58%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "0.0218"
Set-TextValue $ws.Cells.Item(45, 5) "  +2.40%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "8.11"
Set-TextValue $ws.Cells.Item(46, 5) "  +10.46%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "16.87"
Set-TextValue $ws.Cells.Item(47, 5) "  +10.18%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "95.39"
Set-TextValue $ws.Cells.Item(48, 5) "  +4.48%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "1.437.06"
Set-TextValue $ws.Cells.Item(49, 5) "  +4.24%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "2.94"
Set-TextValue $ws.Cells.Item(50, 5) "  +1.34%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "47.42"
Set-TextValue $ws.Cells.Item(51, 5) "  +4.92%  "
